$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 118.52941
$ws.Range("I9").Value = 88.5
$ws.Range("J9").Value = 161.42857
$ws.Range("K9").Value = 88.5
$ws.Range("L9").Value = 161.42857
$ws.Range("M9").Value = 80.5
$ws.Range("N9").Value = -499.42857
$ws.Range("H41").Value = 115.75
$ws.Range("I41").Value = 96.666664
$ws.Range("J41").Value = 134.83333
$ws.Range("K41").Value = 96.666664
$ws.Range("L41").Value = 134.83333
$ws.Range("M41").Value = 343.333336
$ws.Range("N41").Value = -1014.83333
$ws.Range("H53").Value = 43835.695
$ws.Range("I53").Value = 100244.1
$ws.Range("J53").Value = 444.6154
$ws.Range("K53").Value = 100244.1
$ws.Range("L53").Value = 444.6154
$ws.Range("M53").Value = -99607.10000000001
$ws.Range("N53").Value = -1718.6154
$ws.Range("H86").Value = 57126.055
$ws.Range("I86").Value = 92009.91
$ws.Range("J86").Value = 2308.5715
$ws.Range("K86").Value = 92009.91
$ws.Range("L86").Value = 2308.5715
$ws.Range("M86").Value = -90886.91
$ws.Range("N86").Value = -4554.5715
$ws.Range("H89").Value = 57126.055
$ws.Range("I89").Value = 92009.91
$ws.Range("J89").Value = 2308.5715
$ws.Range("K89").Value = 460049.55
$ws.Range("L89").Value = 11542.8575
$ws.Range("M89").Value = -454433.55
$ws.Range("N89").Value = -22774.8575
$ws.Range("H137").Value = 1809.6129
$ws.Range("I137").Value = 1508.4546
$ws.Range("J137").Value = 2545.7778
$ws.Range("K137").Value = 4525.3638
$ws.Range("L137").Value = 7637.3334
$ws.Range("M137").Value = -1975.3638
$ws.Range("N137").Value = -12737.3334
$ws.Range("H138").Value = 1693.0128
$ws.Range("I138").Value = 1264.1578
$ws.Range("J138").Value = 2100.425
$ws.Range("K138").Value = 3792.4734
$ws.Range("L138").Value = 6301.275000000001
$ws.Range("M138").Value = 1347.5266
$ws.Range("N138").Value = -16581.275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17823
$ws.Range("I32").Value = 19212.182
$ws.Range("J32").Value = 6362.25
$ws.Range("K32").Value = 19212.182
$ws.Range("L32").Value = 6362.25
$ws.Range("M32").Value = -18925.182
$ws.Range("N32").Value = -6936.25
$ws.Range("H61").Value = 1894.3846
$ws.Range("I61").Value = 1191.2
$ws.Range("J61").Value = 2333.875
$ws.Range("K61").Value = 1191.2
$ws.Range("L61").Value = 2333.875
$ws.Range("M61").Value = -979.2
$ws.Range("N61").Value = -2757.875
$ws.Range("H74").Value = 1243.9584
$ws.Range("I74").Value = 1202.6316
$ws.Range("J74").Value = 1401
$ws.Range("K74").Value = 1202.6316
$ws.Range("L74").Value = 1401
$ws.Range("M74").Value = -328.6315999999999
$ws.Range("N74").Value = -3149
$ws.Range("H77").Value = 1243.9584
$ws.Range("I77").Value = 1202.6316
$ws.Range("J77").Value = 1401
$ws.Range("K77").Value = 6013.157999999999
$ws.Range("L77").Value = 7005
$ws.Range("M77").Value = -1645.157999999999
$ws.Range("N77").Value = -15741
$ws.Range("H97").Value = 551.2174
$ws.Range("I97").Value = 535.9
$ws.Range("J97").Value = 653.3333
$ws.Range("K97").Value = 535.9
$ws.Range("L97").Value = 653.3333
$ws.Range("M97").Value = -39.89999999999998
$ws.Range("N97").Value = -1645.3333
$ws.Range("H132").Value = 4397.7954
$ws.Range("I132").Value = 5685.731
$ws.Range("J132").Value = 2537.4443
$ws.Range("K132").Value = 17057.193
$ws.Range("L132").Value = 7612.3329
$ws.Range("M132").Value = -14527.193
$ws.Range("N132").Value = -12672.3329
$ws.Range("H135").Value = 30792.363
$ws.Range("J135").Value = 30792.363
$ws.Range("L135").Value = 30792.363
$ws.Range("N135").Value = -40932.363
$ws.Range("H136").Value = 1894.3846
$ws.Range("I136").Value = 1191.2
$ws.Range("J136").Value = 2333.875
$ws.Range("K136").Value = 3573.6
$ws.Range("L136").Value = 7001.625
$ws.Range("M136").Value = -1023.6
$ws.Range("N136").Value = -12101.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2854.6128
$ws.Range("I20").Value = 3130.889
$ws.Range("J20").Value = 2472.077
$ws.Range("K20").Value = 3130.889
$ws.Range("L20").Value = 2472.077
$ws.Range("M20").Value = -2883.889
$ws.Range("N20").Value = -2966.077
$ws.Range("H107").Value = 1070.3334
$ws.Range("I107").Value = 1070.3334
$ws.Range("K107").Value = 1070.3334
$ws.Range("M107").Value = 849.6666
$ws.Range("H134").Value = 41250.848
$ws.Range("I134").Value = 65245.5
$ws.Range("J134").Value = 2859.4
$ws.Range("K134").Value = 195736.5
$ws.Range("L134").Value = 8578.200000000001
$ws.Range("M134").Value = -193201.5
$ws.Range("N134").Value = -13648.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7578656.5
$ws.Range("I31").Value = 2054.7646
$ws.Range("K31").Value = 2054.7646
$ws.Range("M31").Value = -1759.7646
$ws.Range("H34").Value = 7578656.5
$ws.Range("I34").Value = 2054.7646
$ws.Range("K34").Value = 2054.7646
$ws.Range("M34").Value = -1852.7646
$ws.Range("H86").Value = 100001430
$ws.Range("J86").Value = 1552.75
$ws.Range("L86").Value = 1552.75
$ws.Range("N86").Value = -3798.75
$ws.Range("H89").Value = 100001430
$ws.Range("J89").Value = 1552.75
$ws.Range("L89").Value = 7763.75
$ws.Range("N89").Value = -18995.75
$ws.Range("H99").Value = 39486.48
$ws.Range("I99").Value = 85351.914
$ws.Range("J99").Value = 2794.1333
$ws.Range("K99").Value = 85351.914
$ws.Range("L99").Value = 2794.1333
$ws.Range("M99").Value = -83853.914
$ws.Range("N99").Value = -5790.1333
$ws.Range("H126").Value = 39486.48
$ws.Range("I126").Value = 85351.914
$ws.Range("J126").Value = 2794.1333
$ws.Range("K126").Value = 256055.742
$ws.Range("L126").Value = 8382.3999
$ws.Range("M126").Value = -253585.742
$ws.Range("N126").Value = -13322.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4580.625
$ws.Range("I80").Value = 2831.8
$ws.Range("J80").Value = 7495.3335
$ws.Range("K80").Value = 2831.8
$ws.Range("L80").Value = 7495.3335
$ws.Range("M80").Value = -1833.8
$ws.Range("N80").Value = -9491.333500000001
$ws.Range("H83").Value = 4580.625
$ws.Range("I83").Value = 2831.8
$ws.Range("J83").Value = 7495.3335
$ws.Range("K83").Value = 14159
$ws.Range("L83").Value = 37476.6675
$ws.Range("M83").Value = -9167
$ws.Range("N83").Value = -47460.6675
$ws.Range("H97").Value = 852.0714
$ws.Range("H132").Value = 33697.594
$ws.Range("I132").Value = 43947.625
$ws.Range("J132").Value = 2947.5
$ws.Range("K132").Value = 131842.875
$ws.Range("L132").Value = 8842.5
$ws.Range("M132").Value = -129312.875
$ws.Range("N132").Value = -13902.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1514.6666
$ws.Range("I82").Value = 1751.8
$ws.Range("K82").Value = 1751.8
$ws.Range("M82").Value = -1390.8
$ws.Range("H85").Value = 1514.6666
$ws.Range("I85").Value = 1751.8
$ws.Range("K85").Value = 1751.8
$ws.Range("M85").Value = -503.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1059.9
$ws.Range("I132").Value = 1036.6342
$ws.Range("K132").Value = 3109.9026
$ws.Range("M132").Value = -579.9025999999999
